# Auto update Excel log
# Appends new sensor-log rows to the PIR, Proximity and mmWave sheets.
#
# All cells in these log sheets are plain text (t="inlineStr" in the source
# file). Excel's COM layer will happily auto-coerce a literal like
# "2026-01-30" into a real date serial the moment it is poured into
# Range.Value, so the Date column is populated via a temporary text-formula
# ("=""2026-01-30""") and then flattened back to a literal value with
# Copy/PasteSpecial(xlPasteValues) — that round-trip yields a plain string
# cell with no numeric formatting, matching the rest of the sheet.

$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[][]]$Rows   # each inner array: Date, Timestamp, Hour, Location, Value, Status
    )

    # NOTE: this runtime's PowerShell engine does not bind named (-Param)
    # arguments reliably, so this function is always invoked positionally.

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1

    # --- Column A (Date) -------------------------------------------------
    # Write as a text formula first so Excel can't reinterpret the literal
    # "YYYY-MM-DD" string as a date serial, then bake it down to a plain
    # value via paste-special.
    $dateRange = $ws.Range("A$StartRow`:A$endRow")
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $dateValue = $Rows[$i][0]
        $ws.Range("A$r").Formula = '="' + $dateValue + '"'
    }
    $dateRange.Copy()
    $dateRange.PasteSpecial($xlPasteValues)

    # --- Columns B-F -------------------------------------------------------
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $ws.Range("B$r").Value = $Rows[$i][1]
        $ws.Range("C$r").Value = $Rows[$i][2]
        $ws.Range("D$r").Value = $Rows[$i][3]
        $ws.Range("E$r").Value = $Rows[$i][4]
        $ws.Range("F$r").Value = $Rows[$i][5]
    }
}

# ----------------------------------------------------------------------
# PIR sheet: rows 97-113 (Bathroom / No Motion / Inactive)
# ----------------------------------------------------------------------
$pirTimestamps = @(
    "13:06:54",
    "13:06:54",
    "13:06:59",
    "13:07:04",
    "13:07:09",
    "13:07:44",
    "13:07:49",
    "13:07:54",
    "13:07:59",
    "13:08:04",
    "13:08:09",
    "13:08:14",
    "13:08:19",
    "13:08:24",
    "13:08:30",
    "13:08:35",
    "13:08:40"
)

$pirRows = @()
foreach ($ts in $pirTimestamps) {
    $pirRows += , @("2026-01-30", $ts, "13:00", "Bathroom", "No Motion", "Inactive")
}

Add-LogRows "PIR" 97 $pirRows

# ----------------------------------------------------------------------
# Proximity sheet: rows 43-48 (Bathroom Door ENTER/EXIT events)
# ----------------------------------------------------------------------
$proximityRows = @(
    , @("2026-01-30", "13:06:56", "13:00", "Bathroom Door", "EXIT",  "User EXITED Bathroom")
    , @("2026-01-30", "13:07:02", "13:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom")
    , @("2026-01-30", "13:07:45", "13:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom")
    , @("2026-01-30", "13:08:00", "13:00", "Bathroom Door", "EXIT",  "User EXITED Bathroom")
    , @("2026-01-30", "13:08:05", "13:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom")
    , @("2026-01-30", "13:08:20", "13:00", "Bathroom Door", "EXIT",  "User EXITED Bathroom")
)

Add-LogRows "Proximity" 43 $proximityRows

# ----------------------------------------------------------------------
# mmWave sheet: rows 39-40 (Living Room FALL_DETECTED / EMERGENCY)
# ----------------------------------------------------------------------
$mmWaveRows = @(
    , @("2026-01-30", "13:07:43", "13:00", "Living Room", "FALL_DETECTED", "EMERGENCY")
    , @("2026-01-30", "13:07:44", "13:00", "Living Room", "FALL_DETECTED", "EMERGENCY")
)

Add-LogRows "mmWave" 39 $mmWaveRows
